$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.748.85"
$ws.Range("E2").Value = "  +1.36%  "
$ws.Range("D3").Value = "2.496.83"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'587.17"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("D6").Value = "'176.74"
$ws.Range("E6").Value = "  +5.01%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.516"
$ws.Range("E8").Value = "  +1.58%  "
$ws.Range("D9").Value = "'0.140"
$ws.Range("E9").Value = "  +5.48%  "
$ws.Range("E10").Value = "  +0.80%  "
$ws.Range("D11").Value = "'0.339"
$ws.Range("E11").Value = "  +4.18%  "
$ws.Range("D12").Value = "'4.95"
$ws.Range("E12").Value = "  +1.84%  "
$ws.Range("D13").Value = "'25.74"
$ws.Range("E13").Value = "  +2.89%  "
$ws.Range("E14").Value = "  +1.01%  "
$ws.Range("D15").Value = "67.536.50"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("E16").Value = "  +2.98%  "
$ws.Range("D17").Value = "2.488.43"
$ws.Range("E17").Value = "  +2.77%  "
$ws.Range("D18").Value = "'11.08"
$ws.Range("E18").Value = "  +1.98%  "
$ws.Range("D19").Value = "'7.44"
$ws.Range("E19").Value = "  +1.68%  "
$ws.Range("D20").Value = "'352.05"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").Value = "'4.08"
$ws.Range("E21").Value = "  +2.36%  "
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "'70.77"
$ws.Range("E23").Value = "  +3.51%  "
$ws.Range("D24").Value = "'4.24"
$ws.Range("E24").Value = "  +2.25%  "
$ws.Range("D25").Value = "'1.79"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("D26").Value = "'9.21"
$ws.Range("E26").Value = "  +2.58%  "
$ws.Range("D27").Value = "2.623.47"
$ws.Range("E27").Value = "  +1.95%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("D29").Value = "0.0₃0914"
$ws.Range("E29").Value = "  +3.57%  "
$ws.Range("D30").Value = "'513.71"
$ws.Range("E30").Value = "  +1.85%  "
$ws.Range("D31").Value = "'7.87"
$ws.Range("E31").Value = "  +4.35%  "
$ws.Range("E32").Value = "  +3.86%  "
$ws.Range("E33").Value = "  +2.13%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  +8.72%  "
$ws.Range("D36").Value = "'161.31"
$ws.Range("E36").Value = "  +2.07%  "
$ws.Range("D37").Value = "'18.71"
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("D38").Value = "'18.41"
$ws.Range("E38").Value = "  +1.94%  "
$ws.Range("E39").Value = "  +2.29%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "'1.74"
$ws.Range("E41").Value = "  +5.34%  "
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").Value = "'0.330"
$ws.Range("E42").Value = "  +3.01%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").Value = "'4.88"
$ws.Range("E43").Value = "  +3.73%  "
$ws.Range("D44").Value = "'2.44"
$ws.Range("E44").Value = "  +4.85%  "
$ws.Range("D45").Value = "'144.39"
$ws.Range("E45").Value = "  +3.26%  "
$ws.Range("D46").Value = "'3.52"
$ws.Range("E46").Value = "  +3.65%  "
$ws.Range("E47").Value = "  +4.73%  "
$ws.Range("D48").Value = "'0.515"
$ws.Range("E48").Value = "  +2.64%  "
$ws.Range("D49").Value = "'0.0745"
$ws.Range("E49").Value = "  +2.63%  "
$ws.Range("D50").Value = "'1.59"
$ws.Range("E50").Value = "  +2.65%  "
$ws.Range("D51").Value = "'0.587"
$ws.Range("E51").Value = "  +2.12%  "
